$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins / Losses / Ties) right after the
# existing "Unnamed: 28" column (AC), in columns AD, AE, AF of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / centered / bordered header formatting used by the rest
# of row 1 (e.g. A1) by copying its format onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the team record (Wins=63, Losses=97, Ties=1) for every player
# row in the sheet (rows 2 through 39).
$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 63
    $ws.Cells.Item($r, 31).Value = 97
    $ws.Cells.Item($r, 32).Value = 1
}
